$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Merge the "February " + "3, 2023" runs into a single run.
#    Both runs already share identical run formatting, and the text
#    contains no apostrophes/quotes, so a straightforward Find &
#    Replace safely merges them while preserving the shared rPr.
# ------------------------------------------------------------------
$d.Content.Find.Execute("February 3, 2023", $true, $false, $false, $false, $false, $true, 1, $false, "February 3, 2023", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Merge "Have you ever seen so much pr" + "e" + "-refutation as in
#    Ed Kozlowski's recent letter? I" into one run. This text has a
#    straight apostrophe that must be preserved (not turned into a
#    curly quote), so instead of Find & Replace (which smart-quotes
#    typed text) we delete the paragraph's text and retype it via
#    Range.InsertAfter, then nudge Bold on/off (while excluding the
#    paragraph mark) purely to coax the engine into emitting the
#    empty <w:rPr/> element that was present in the original XML.
# ------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Have you ever seen so much pre-refutation as in Ed Kozlowski's recent letter? I") {
        $idx = $i
        break
    }
}
$p = $d.Paragraphs.Item($idx)
$start = $p.Range.Start
$end = $p.Range.End
$delRange = $d.Range($start, $end - 1)
$delRange.Delete()
$p2 = $d.Paragraphs.Item($idx)
$insPoint = $d.Range($p2.Range.Start, $p2.Range.Start)
$insPoint.InsertAfter("Have you ever seen so much pre-refutation as in Ed Kozlowski's recent letter? I")
$p3 = $d.Paragraphs.Item($idx)
$fmtRange = $p3.Range
$fmtRange.MoveEnd(1, -1) | Out-Null
$fmtRange.Font.Bold = 1
$p4 = $d.Paragraphs.Item($idx)
$fmtRange2 = $p4.Range
$fmtRange2.MoveEnd(1, -1) | Out-Null
$fmtRange2.Font.Bold = 0

# ------------------------------------------------------------------
# 3. Split "(860) 335-2869 (C)" into "(860) 335-2869 " and "(C)" runs.
#    Toggling Bold on the "(C)" sub-range and back off forces the
#    engine to split the run at that boundary while leaving both
#    runs with an empty (but present) <w:rPr/>.
# ------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "(860) 335-2869 (C)") {
        $idx2 = $i
        break
    }
}
$pPhone = $d.Paragraphs.Item($idx2)
$phoneStart = $pPhone.Range.Start
$cRange = $d.Range($phoneStart + 15, $phoneStart + 18)
$cRange.Font.Bold = 1
$cRange2 = $d.Range($phoneStart + 15, $phoneStart + 18)
$cRange2.Font.Bold = 0

# ------------------------------------------------------------------
# 4. Add two new paragraphs after the phone-number paragraph:
#      a) an empty paragraph (same indentation as the line above)
#      b) a right-aligned paragraph with the co-chair note
# ------------------------------------------------------------------
$pPhone = $d.Paragraphs.Item($idx2)
$afterPhone = $pPhone.Range
$afterPhone.InsertParagraphAfter() | Out-Null

$paras = $d.Paragraphs
$blankPara = $paras.Item($paras.Count)
$blankPara.Range.InsertParagraphAfter() | Out-Null

$paras = $d.Paragraphs
$notePara = $paras.Item($paras.Count)
$notePara.Format.LeftIndent = 0
$notePara.Format.Alignment = 2
$notePara.Range.Text = "(Co-chair of Andover DTC and member of RHAM BOE, but writing as an individual.) "
